$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.094902791605733228,
    -0.0099999995071371472,
    -0.008999999512466772,
    0.28399241863110802,
    -0.0059999995260593408,
    -0.0059999995091537528,
    -0.01999999942842301,
    -0.01999999942467845,
    -0.0059999995006059237,
    0.030061936047211191,
    -0.0044999995057182218,
    -0.0059999994953967573,
    -0.0059999994875532536,
    -0.011999999451324683,
    -0.0059999994837092174,
    -0.0059999994821566816,
    -0.0059999994800810086,
    -0.0089999994626328572,
    -0.0089999995116145648,
    -0.008999999498438882,
    -0.0089999994958862573,
    -0.057853776727209372,
    -0.0089999995083722695,
    -0.041999999313302894,
    -0.04199999930976972,
    -0.0059999995075408208,
    -0.0059999995053505728,
    -0.0059999994955841629,
    -0.010135424125259362,
    -0.021335050439405645,
    -0.014999999430367339,
    -0.020999999395599822,
    -0.0059999994805099988
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
